# Updated Array and linkedlist pages
$wb = $excel.ActiveWorkbook

# --- Adjust existing "Text_Try_Editor" sheet view (it is no longer the active tab) ---
$wsEditor = $wb.Worksheets.Item("Text_Try_Editor")
$wsEditor.Activate() | Out-Null
$wsEditor.Range("F12").Select() | Out-Null

# --- Add the new "practiceQuestions" sheet after the last existing sheet ---
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "practiceQuestions"

$codeSearch = @'
my_list = [1, 2, 3, 4, 5] # Element to check
element = 3
# Check if the element is in the list
if element in my_list:
    print("The element exists in the list.")
else:
    print("The element does not exist in the list.")
'@

$codeEven = @'
example = [22, 234, 2463]
def is_even(value):
    length = len(str(value))
    return length % 2 == 0
count = 0
for i in example:
    if is_even(i):
        count += 1
print(count)
'@

# Row 1 - headers (shared-string order matches the authored file: invalidText, validText, Links)
$ws.Range("A1").Value = "invalidText"
$ws.Range("B1").Value = "validText"
$ws.Range("C1").Value = "Links"

# Column C text for rows 2-5 first (keeps shared-string insertion order identical to the source workbook)
$ws.Range("C2").Value = "Search the array"
$ws.Range("C3").Value = "Max Consecutive Ones"
$ws.Range("C4").Value = "Find Numbers with Even Number"
$ws.Range("C5").Value = "Squares of a Sorted Array"

# Column B code snippets (row 2 gets its own snippet, rows 3-5 share one snippet)
$ws.Range("B2").Value = $codeSearch
$ws.Range("B3").Value = $codeEven
$ws.Range("B4").Value = $codeEven
$ws.Range("B5").Value = $codeEven

# Column A repeats the "hello" shared string already present in the workbook
$ws.Range("A2").Value = "hello"
$ws.Range("A3").Value = "hello"
$ws.Range("A4").Value = "hello"
$ws.Range("A5").Value = "hello"

# Wrap text + row heights for the code cells
$ws.Range("B2:B5").WrapText = $true
$ws.Rows.Item(2).RowHeight = 101.5
$ws.Rows.Item(3).RowHeight = 174
$ws.Rows.Item(4).RowHeight = 174
$ws.Rows.Item(5).RowHeight = 174

# --- column widths (characters) ---
$ws.Columns.Item(1).ColumnWidth = 9
$ws.Columns.Item(2).ColumnWidth = 59.666666666666664
$ws.Columns.Item(3).ColumnWidth = 27.833333333333332

# --- view state: select B5 and make this the active/displayed tab ---
$ws.Activate() | Out-Null
$ws.Range("B5").Select() | Out-Null
